$d = $word.ActiveDocument

# 1. Update the cached TIME field text (date string).
$d.Content.Find.Execute("lunes, 28 de agosto de 2017", $true, $false, $false, $false, $false,
                         $true, 1, $false, "miércoles, 27 de septiembre de 2017", 2)

# 2. Update the cached MERGEFIELD "APELLIDO" display text.
$d.Content.Find.Execute("Pez", $true, $false, $false, $false, $false,
                         $true, 1, $false, "«APELLIDO»", 2)

# 3. Update the cached MERGEFIELD "NOMBRE" display text.
$d.Content.Find.Execute("Andrea", $true, $false, $false, $false, $false,
                         $true, 1, $false, "«NOMBRE»", 2)

# 4. Move the "_GoBack" bookmark from the middle of the "La Dirección..."
#    paragraph down to the start of the "Dirección de Pedagogía
#    Universitaria" signature paragraph, and merge the two runs that used
#    to straddle it into a single run (re-typing the whole sentence causes
#    Word to coalesce the run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Content.Find.Execute("La Dirección de Pedagogía Universitaria de la Universidad Nacional de la Matanza informa que se encuentra abierta la inscripción a las asignaturas para el cuatrimestre en curso.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "La Dirección de Pedagogía Universitaria de la Universidad Nacional de la Matanza informa que se encuentra abierta la inscripción a las asignaturas para el cuatrimestre en curso.",
                         2)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text
    if (($text -like "Dirección de Pedagogía Universitaria*") -and -not ($text -like "La Dirección*")) {
        $target = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $target)
        break
    }
}
